$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spectral-type (ST) corrections:
#   Cen X-3 (row 6) ST changed from "O9III" to "O6.5III"
#   Vela X-1 (row 5) ST changed from "B0.5Ia" to "B0.5Ib"
# Written in this order so the new shared-string entries land in the same
# sequence as the target workbook (O6.5III first, then B0.5Ib).
$ws.Range("B6").Value = "O6.5III"
$ws.Range("B5").Value = "B0.5Ib"

# Update the sheet's active selection/cursor position
$ws.Range("B20").Select()
